{"js": "// The document has a handful of label/value lines built from multiple\n// <w:t> runs joined by <w:br/>. Several of those runs carry a stray\n// trailing two-space sequence (an old formatting leftover) that needs to\n// be trimmed, and the signature block opens with two consecutive line\n// breaks where only one is wanted.\n\nconst body = context.document.body;\n\n// Exact (pre -> post) text replacements: trim the trailing double space\n// from each of these run texts.\nconst trims = [\n  \"Sponsor: [Sponsor Name]  \",\n  \"Contact Email: [Email]  \",\n  \"Asset Sponsored: [Asset Type]  \",\n  \"Location: [Location Scope]  \",\n  \"Tier: [Tier Level]  \",\n  \"Impressions Expected: [Impressions]  \",\n];\n\nfor (const original of trims) {\n  const trimmed = original.replace(/\\s+$/, \"\");\n  const found = body.search(original, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const hit of found.items) {\n    hit.insertText(trimmed, \"Replace\");\n  }\n  await context.sync();\n}\n\n// Collapse the double line-break (\"\\u000b\\u000b\") right before the\n// \"Sponsor Signature\" text down to a single line break.\nconst doubleBreak = body.search(\"\\u000b\\u000bSponsor Signature\", { matchCase: true });\ndoubleBreak.load(\"items\");\nawait context.sync();\n\nfor (const hit of doubleBreak.items) {\n  hit.insertText(\"\\u000bSponsor Signature\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The document has a handful of label/value lines built from multiple\n# runs joined by manual line breaks. Several of those runs carry a stray\n# trailing two-space sequence (an old formatting leftover) that needs to\n# be trimmed, and the signature block opens with two consecutive line\n# breaks where only one is wanted.\n\n$d = $word.ActiveDocument\n\n# Manual line break character, as it appears in Range.Text (same glyph\n# that <w:br/> round-trips through, vertical tab / Chr(11)).\n$vtab = [char]11\n\nfunction Replace-ExactText($searchText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n# Trim the trailing double space from each of these run texts.\nReplace-ExactText \"Sponsor: [Sponsor Name]  \" \"Sponsor: [Sponsor Name]\"\nReplace-ExactText \"Contact Email: [Email]  \" \"Contact Email: [Email]\"\nReplace-ExactText \"Asset Sponsored: [Asset Type]  \" \"Asset Sponsored: [Asset Type]\"\nReplace-ExactText \"Location: [Location Scope]  \" \"Location: [Location Scope]\"\nReplace-ExactText \"Tier: [Tier Level]  \" \"Tier: [Tier Level]\"\nReplace-ExactText \"Impressions Expected: [Impressions]  \" \"Impressions Expected: [Impressions]\"\n\n# Collapse the double line-break right before \"Sponsor Signature\" down to\n# a single line break.\nReplace-ExactText \"${vtab}${vtab}Sponsor Signature\" \"${vtab}Sponsor Signature\"\n"}
